$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (K), rows 2-37, replacing the previous Strike# derived values
$gValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 2
    6  = 1
    7  = 1
    8  = 0
    9  = 1
    10 = 2
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 1
    18 = 2
    19 = 0
    20 = 2
    21 = 0
    22 = 2
    23 = 1
    24 = 1
    25 = 2
    26 = 0
    27 = 0
    28 = 2
    29 = 1
    30 = 2
    31 = 1
    32 = 3
    33 = 2
    34 = 0
    35 = 0
    36 = 1
    37 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
